$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 ---
$ws.Range("A1").Value = 158062
$ws.Range("B1").Value = "Assisti a uma vídeo-aula no Youtube"
$ws.Range("C1").Value = "https://www.youtube.com/asd"

# --- Row 2 ---
$ws.Range("A2").Value = 159070
$ws.Range("B2").Value = "Vídeoaula1_atualizado"
$ws.Range("C2").Value = "https://www.youtube.com/watch?v=PKMm-cHe56g "

# --- Row 3 (new) ---
$ws.Range("A3").Value = 162079
$ws.Range("B3").Value = "Vídeoaula10_atualizado"
$ws.Range("C3").Value = "https://www.youtube.com/watch?v=PamJA8e56g "

# New row 3 formatting: B3 gets wrap text with general horizontal alignment
$ws.Range("B3").WrapText = $true
$ws.Range("B3").HorizontalAlignment = 1

# --- Hyperlink: move from C2 to C1, repointing to the new URL ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C1"), "https://www.youtube.com/asd", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://www.youtube.com/asd")

# Restore C1's original (non-hyperlink) look by re-applying C2's format,
# since Excel auto-applies the "Hyperlink" style when a hyperlink is added.
$ws.Range("C2").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0 | Out-Null

# --- Selection / active cell matches the authored file ---
$ws.Range("C3").Select() | Out-Null
